$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09769597576622248
$ws.Range("H2").Value = -25.66167679853846
$ws.Range("I2").Value = 33.29302245774368
$ws.Range("G3").Value = 0.09605050813595715
$ws.Range("H3").Value = 7.87683365147724
$ws.Range("G4").Value = -0.5892784171921228
$ws.Range("H4").Value = 4.235848991760185
$ws.Range("G5").Value = -0.6458596125124123
$ws.Range("H5").Value = -5.775754653012443
$ws.Range("G6").Value = 0.1794482936002075
$ws.Range("H6").Value = -27.08548627260313
$ws.Range("G7").Value = 0.3486035798085087
$ws.Range("H7").Value = 112.7941922549925
$ws.Range("G8").Value = 0.1246453306811028
$ws.Range("H8").Value = -24.59854897190949
$ws.Range("G9").Value = 0.1709887539645877
$ws.Range("H9").Value = -12.35505215578009
$ws.Range("G10").Value = -0.0895536722284744
$ws.Range("H10").Value = -56.70935403281615
$ws.Range("G11").Value = -0.0910819325612508
$ws.Range("H11").Value = 23.31125821592646
$ws.Range("G12").Value = 0.2186243686942796
$ws.Range("H12").Value = 37.4742819240659
$ws.Range("G13").Value = 0.2454527102231195
$ws.Range("H13").Value = 19.34777052390841
$ws.Range("G14").Value = 0.1633314461548659
$ws.Range("H14").Value = -13.75149554692259
$ws.Range("G15").Value = 0.2186041130327111
$ws.Range("H15").Value = -12.52526051017337
$ws.Range("G16").Value = 0.0424890268819675
$ws.Range("H16").Value = 16.47132489198356
$ws.Range("G17").Value = 0.0498025332165964
$ws.Range("H17").Value = 40.40514451389804
$ws.Range("G18").Value = 0.1459611505117161
$ws.Range("H18").Value = -15.78235896213192
$ws.Range("G19").Value = 0.1974803338811424
$ws.Range("H19").Value = 57.03833868897304
$ws.Range("G20").Value = 0.1005002179685333
$ws.Range("H20").Value = -12.34324978037785
$ws.Range("G21").Value = 0.1083073201965206
$ws.Range("H21").Value = 7.87707842307011
$ws.Range("G22").Value = 0.06356757921470642
$ws.Range("H22").Value = -32.51643676584996
$ws.Range("G23").Value = 0.05412745353925746
$ws.Range("H23").Value = -50.10834215829632
$ws.Range("G24").Value = -0.2317487216623107
$ws.Range("H24").Value = -85.95782166206843
$ws.Range("G25").Value = -0.2398536863299895
$ws.Range("H25").Value = -7.821658780457375
$ws.Range("G26").Value = 0.1722455443511939
$ws.Range("H26").Value = 8.342920487110876
$ws.Range("G27").Value = 0.2082760644615324
$ws.Range("H27").Value = 3.895780302769106
$ws.Range("G28").Value = 0.03661646163937396
$ws.Range("H28").Value = 555.9702218303879
$ws.Range("G29").Value = 0.07998696664967388
$ws.Range("H29").Value = 420.1561311674374
